$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 133
$ws.Cells.Item(4, 9).Value = 8.333333
$ws.Cells.Item(4, 10).Value = 320
$ws.Cells.Item(4, 11).Value = 8.333333
$ws.Cells.Item(4, 12).Value = 320
$ws.Cells.Item(4, 13).Value = 105.666667
$ws.Cells.Item(4, 14).Value = -548

$ws.Cells.Item(39, 8).Value = 542.3333
$ws.Cells.Item(39, 9).Value = 485.125
$ws.Cells.Item(39, 10).Value = 1000
$ws.Cells.Item(39, 11).Value = 1455.375
$ws.Cells.Item(39, 12).Value = 3000
$ws.Cells.Item(39, 13).Value = -1159.375
$ws.Cells.Item(39, 14).Value = -3592

$ws.Cells.Item(64, 8).Value = 5999.125
$ws.Cells.Item(64, 9).Value = 4999.3335
$ws.Cells.Item(64, 11).Value = 4999.3335
$ws.Cells.Item(64, 13).Value = -4751.3335

$ws.Cells.Item(67, 8).Value = 5999.125
$ws.Cells.Item(67, 9).Value = 4999.3335
$ws.Cells.Item(67, 11).Value = 4999.3335
$ws.Cells.Item(67, 13).Value = -4141.3335

$ws.Cells.Item(74, 8).Value = 7145.7144
$ws.Cells.Item(74, 9).Value = 4956.75
$ws.Cells.Item(74, 11).Value = 4956.75
$ws.Cells.Item(74, 13).Value = -4020.75

$ws.Cells.Item(77, 8).Value = 7145.7144
$ws.Cells.Item(77, 9).Value = 4956.75
$ws.Cells.Item(77, 11).Value = 24783.75
$ws.Cells.Item(77, 13).Value = -20103.75

$ws.Cells.Item(98, 8).Value = 5057.0557
$ws.Cells.Item(98, 9).Value = 2955.8333
$ws.Cells.Item(98, 10).Value = 9259.5
$ws.Cells.Item(98, 11).Value = 2955.8333
$ws.Cells.Item(98, 12).Value = 9259.5
$ws.Cells.Item(98, 13).Value = -1457.8333
$ws.Cells.Item(98, 14).Value = -12255.5

$ws.Cells.Item(122, 8).Value = 5057.0557
$ws.Cells.Item(122, 9).Value = 2955.8333
$ws.Cells.Item(122, 10).Value = 9259.5
$ws.Cells.Item(122, 11).Value = 8867.499899999999
$ws.Cells.Item(122, 12).Value = 27778.5
$ws.Cells.Item(122, 13).Value = -6417.499899999999
$ws.Cells.Item(122, 14).Value = -32678.5

$ws.Cells.Item(132, 8).Value = 24398162
$ws.Cells.Item(132, 9).Value = 27030912
$ws.Cells.Item(132, 11).Value = 81092736
$ws.Cells.Item(132, 13).Value = -81090206

$ws.Cells.Item(137, 8).Value = 2213.4614
$ws.Cells.Item(137, 9).Value = 2197.7896
$ws.Cells.Item(137, 11).Value = 6593.3688
$ws.Cells.Item(137, 13).Value = -4043.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7714.911
$ws.Cells.Item(32, 9).Value = 7722.436
$ws.Cells.Item(32, 11).Value = 7722.436
$ws.Cells.Item(32, 13).Value = -7435.436

$ws.Cells.Item(45, 8).Value = 2137.1538
$ws.Cells.Item(45, 10).Value = 3000
$ws.Cells.Item(45, 12).Value = 3000
$ws.Cells.Item(45, 14).Value = -3754

$ws.Cells.Item(110, 8).Value = 1785.7826
$ws.Cells.Item(110, 9).Value = 1609.2106
$ws.Cells.Item(110, 11).Value = 1609.2106
$ws.Cells.Item(110, 13).Value = 435.7893999999999

$ws.Cells.Item(122, 8).Value = 3650
$ws.Cells.Item(122, 9).Value = 3650
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 10950
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -8500
$ws.Cells.Item(122, 14).Value = $null

$ws.Cells.Item(132, 8).Value = 7245.4614
$ws.Cells.Item(132, 9).Value = 2242
$ws.Cells.Item(132, 11).Value = 6726
$ws.Cells.Item(132, 13).Value = -4196

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2535.889
$ws.Cells.Item(20, 9).Value = 2585.6428
$ws.Cells.Item(20, 10).Value = 2482.3076
$ws.Cells.Item(20, 11).Value = 2585.6428
$ws.Cells.Item(20, 12).Value = 2482.3076
$ws.Cells.Item(20, 13).Value = -2338.6428
$ws.Cells.Item(20, 14).Value = -2976.3076

$ws.Cells.Item(60, 8).Value = 35250
$ws.Cells.Item(60, 10).Value = 47875
$ws.Cells.Item(60, 12).Value = 47875
$ws.Cells.Item(60, 14).Value = -49073

$ws.Cells.Item(99, 8).Value = 4263.6
$ws.Cells.Item(99, 9).Value = 2507.4
$ws.Cells.Item(99, 11).Value = 2507.4
$ws.Cells.Item(99, 13).Value = -1009.4

$ws.Cells.Item(105, 8).Value = 3575.1365
$ws.Cells.Item(105, 9).Value = 2598.0557
$ws.Cells.Item(105, 11).Value = 2598.0557
$ws.Cells.Item(105, 13).Value = -851.0556999999999

$ws.Cells.Item(107, 8).Value = 808.6896400000001
$ws.Cells.Item(107, 9).Value = 649.5599999999999
$ws.Cells.Item(107, 10).Value = 1803.25
$ws.Cells.Item(107, 11).Value = 649.5599999999999
$ws.Cells.Item(107, 12).Value = 1803.25
$ws.Cells.Item(107, 13).Value = 1270.44
$ws.Cells.Item(107, 14).Value = -5643.25

$ws.Cells.Item(132, 8).Value = 123333
$ws.Cells.Item(132, 10).Value = 123333
$ws.Cells.Item(132, 12).Value = 123333
$ws.Cells.Item(132, 14).Value = -133453

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2497.9375
$ws.Cells.Item(31, 9).Value = 1971.1163
$ws.Cells.Item(31, 11).Value = 1971.1163
$ws.Cells.Item(31, 13).Value = -1676.1163

$ws.Cells.Item(34, 8).Value = 2497.9375
$ws.Cells.Item(34, 9).Value = 1971.1163
$ws.Cells.Item(34, 11).Value = 1971.1163
$ws.Cells.Item(34, 13).Value = -1769.1163

$ws.Cells.Item(51, 8).Value = 71175
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 71175
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 71175
$ws.Cells.Item(51, 13).Value = $null
$ws.Cells.Item(51, 14).Value = -72647

$ws.Cells.Item(61, 8).Value = 71175
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 71175
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 71175
$ws.Cells.Item(61, 13).Value = $null
$ws.Cells.Item(61, 14).Value = -71871

$ws.Cells.Item(105, 8).Value = 1778
$ws.Cells.Item(105, 9).Value = 1487.8
$ws.Cells.Item(105, 11).Value = 1487.8
$ws.Cells.Item(105, 13).Value = 259.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 239.8
$ws.Cells.Item(26, 9).Value = 199
$ws.Cells.Item(26, 10).Value = 250
$ws.Cells.Item(26, 11).Value = 597
$ws.Cells.Item(26, 12).Value = 750
$ws.Cells.Item(26, 13).Value = -309
$ws.Cells.Item(26, 14).Value = -1326

$ws.Cells.Item(39, 8).Value = 5347.5
$ws.Cells.Item(39, 10).Value = 8366
$ws.Cells.Item(39, 12).Value = 25098
$ws.Cells.Item(39, 14).Value = -25686

$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 13).Value = $null

$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 13).Value = $null

$ws.Cells.Item(118, 8).Value = 8333.333000000001
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 13).Value = $null

$ws.Cells.Item(128, 8).Value = 98999.336
$ws.Cells.Item(128, 9).Value = 98999.336
$ws.Cells.Item(128, 11).Value = 296998.008
$ws.Cells.Item(128, 13).Value = -292018.008

$ws.Cells.Item(137, 8).Value = 3584.0435
$ws.Cells.Item(137, 9).Value = 980.2222
$ws.Cells.Item(137, 10).Value = 5257.9287
$ws.Cells.Item(137, 11).Value = 2940.6666
$ws.Cells.Item(137, 12).Value = 15773.7861
$ws.Cells.Item(137, 13).Value = 2159.3334
$ws.Cells.Item(137, 14).Value = -25973.7861

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 697.4857
$ws.Cells.Item(97, 9).Value = 662.04
$ws.Cells.Item(97, 10).Value = 786.1
$ws.Cells.Item(97, 11).Value = 662.04
$ws.Cells.Item(97, 12).Value = 786.1
$ws.Cells.Item(97, 13).Value = -166.04
$ws.Cells.Item(97, 14).Value = -1778.1

$ws.Cells.Item(122, 8).Value = 2424.125
$ws.Cells.Item(122, 9).Value = 2181.3572
$ws.Cells.Item(122, 10).Value = 4123.5
$ws.Cells.Item(122, 11).Value = 6544.071599999999
$ws.Cells.Item(122, 12).Value = 12370.5
$ws.Cells.Item(122, 13).Value = -4094.071599999999
$ws.Cells.Item(122, 14).Value = -17270.5

$ws.Cells.Item(132, 8).Value = 5218.6
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = $null

$ws.Cells.Item(141, 8).Value = 90736
$ws.Cells.Item(141, 10).Value = 90736
$ws.Cells.Item(141, 12).Value = 90736
$ws.Cells.Item(141, 14).Value = -101096

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 132.39131
$ws.Cells.Item(55, 10).Value = 176.27272
$ws.Cells.Item(55, 12).Value = 176.27272
$ws.Cells.Item(55, 14).Value = -522.2727199999999

$ws.Cells.Item(93, 8).Value = 4511.5713
$ws.Cells.Item(93, 9).Value = 4826.2
$ws.Cells.Item(93, 11).Value = 4826.2
$ws.Cells.Item(93, 13).Value = -3578.2

$ws.Cells.Item(122, 8).Value = 2832.111
$ws.Cells.Item(122, 9).Value = 2242.7693
$ws.Cells.Item(122, 10).Value = 4364.4
$ws.Cells.Item(122, 11).Value = 6728.3079
$ws.Cells.Item(122, 12).Value = 13093.2
$ws.Cells.Item(122, 13).Value = -4278.3079
$ws.Cells.Item(122, 14).Value = -17993.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1445.3846
$ws.Cells.Item(96, 9).Value = 793
$ws.Cells.Item(96, 11).Value = 793
$ws.Cells.Item(96, 13).Value = 580

$ws.Cells.Item(100, 8).Value = 1262.9474
$ws.Cells.Item(100, 9).Value = 1232.9032
$ws.Cells.Item(100, 10).Value = 1396
$ws.Cells.Item(100, 11).Value = 2465.8064
$ws.Cells.Item(100, 12).Value = 2792
$ws.Cells.Item(100, 13).Value = -1924.8064
$ws.Cells.Item(100, 14).Value = -3874

$ws.Cells.Item(126, 8).Value = 2623.7778
$ws.Cells.Item(126, 9).Value = 2339.75
$ws.Cells.Item(126, 11).Value = 7019.25
$ws.Cells.Item(126, 13).Value = -4549.25

$ws.Cells.Item(136, 8).Value = 8001.276
$ws.Cells.Item(136, 10).Value = 1949.25
$ws.Cells.Item(136, 12).Value = 5847.75
$ws.Cells.Item(136, 14).Value = -10947.75
